$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lifts")

$ws.Cells.Item(51, 1).Value = 48
$ws.Cells.Item(51, 2).Value = "WAREHOUSE"
$ws.Cells.Item(51, 3).Value = "S"
$ws.Cells.Item(51, 4).Value = "O"
$ws.Cells.Item(51, 5).Value = "Teinipeili selfie piti ottaa pitkästä aikaa"
